$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add number 5 in cell A2 (the "5th cell" following A1..D1 being cells 1-4)
$ws.Range("A2").Value = 5

# Move the active selection to F5, matching the recorded final cursor position
$ws.Range("F5").Select()
